# Update the "想去人数" (number of people interested) counts for two events
# that appear in both the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F3").Value = 109
    $ws.Range("F5").Value = 68
}
